$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the 5 wrong "pasadas funales" (final pass) values ---
$ws.Range("D3").Value = 5
$ws.Range("D7").Value = 17
$ws.Range("C12").Value = 33
$ws.Range("C13").Value = 45
$ws.Range("D20").Value = 45

# --- 2. Recolor the shared green fill's background (fgColor stays FFE2EFDA,
#         bgColor becomes solid black instead of the default indexed color) ---
$ws.Range("A13").Interior.PatternColor = 0

# --- 3. Clean up the formatting: every data cell (A2:E21) should share the
#         same plain style (green fill, black font, General number format) -
#         copy the now-updated reference format from A13 onto the whole block ---
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A2:E21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 4. Update the active selection shown when the sheet is reopened ---
$ws.Range("A2:E21").Select()

# --- 5. Best-effort: remember the workbook window position ---
$win = $wb.Windows.Item(1)
$win.Left = 990
$win.Top = 4365
